# Applies the "fix bug of the interval in changing character" commit.
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Settings sheet
# ---------------------------------------------------------------------------
$settings = $wb.Worksheets.Item("Settings")

$settings.Range("B10").Value = "雷电模拟器"
$settings.Range("B10").HorizontalAlignment = -4108
$settings.Range("B10").VerticalAlignment = -4108

$settings.Range("B6").Value = 0
$settings.Range("B5").Value = "CBA"

$settings.Range("B2").Value = 1.25
$settings.Range("B3").Value = 1000
$settings.Range("B4").Value = 1
$settings.Range("B7").Value = 17
$settings.Range("B11").Value = 1

# ---------------------------------------------------------------------------
# Side1 sheet
# ---------------------------------------------------------------------------
$side1 = $wb.Worksheets.Item("Side1")

$side1.Range("B5").Value = 1
$side1.Range("C5").Value = 0

$side1.Range("B6").Value = 1
$side1.Range("C6").Value = 1

$side1.Range("B7").Value = 1
$side1.Range("C7").Value = 1

$side1.Range("B9").Value = 1
$side1.Range("C9").Value = 0

$side1.Range("B11").Value = 1
$side1.Range("C11").Value = 1

$side1.Range("B14").Value = 1
$side1.Range("C14").Value = 3
$side1.Range("D14").Value = 4

$side1.Range("B17").Value = 1
$side1.Range("C17").Value = 1

$side1.Range("B22").Value = 29

# ---------------------------------------------------------------------------
# Side2 sheet
# ---------------------------------------------------------------------------
$side2 = $wb.Worksheets.Item("Side2")

$side2.Range("B5").Value = ""
$side2.Range("C5").Value = ""

$side2.Range("B8").Value = ""
$side2.Range("C8").Value = ""

$side2.Range("B9").Value = ""
$side2.Range("C9").Value = ""

$side2.Range("B10").Value = ""
$side2.Range("C10").Value = ""

$side2.Range("B11").Value = ""
$side2.Range("C11").Value = ""

$side2.Range("B19").Value = ""
$side2.Range("C19").Value = ""

$side2.Range("B22").Value = 35

# ---------------------------------------------------------------------------
# Side3 sheet
# ---------------------------------------------------------------------------
$side3 = $wb.Worksheets.Item("Side3")

$side3.Range("B6").Value = ""
$side3.Range("C6").Value = ""

$side3.Range("B10").Value = 1
$side3.Range("C10").Value = 0

$side3.Range("B12").Value = 1
$side3.Range("C12").Value = 0

$side3.Range("B13").Value = ""
$side3.Range("C13").Value = ""

$side3.Range("B19").Value = 1
$side3.Range("C19").Value = 0

$side3.Range("B20").Value = ""
$side3.Range("C20").Value = ""

# ---------------------------------------------------------------------------
# Selections / active sheet bookkeeping
# ---------------------------------------------------------------------------
$settings.Range("A6").Select()
$side1.Range("A15").Select()
$side2.Range("B20").Select()
$side3.Range("A23").Select()

$side1.Activate()
